# River trend results update - May 2024
# Updates the recalculated trend statistics (confidence, proportions,
# medians, Sen slopes, confidence intervals, percent annual change) and
# the derived "confidence of improving trend" / analysis-note text for
# each parameter row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Dissolved Oxygen Concentration
$ws.Range("F2").Value = 0.7777805887479819
$ws.Range("H2").Value = 0.946428571428571
$ws.Range("J2").Value = 9.875
$ws.Range("K2").Value = 0.0903090659340658
$ws.Range("L2").Value = -0.116522432242832
$ws.Range("M2").Value = 0.246651172932376
$ws.Range("N2").Value = 0.914522186674084

# Row 3 - Dissolved Reactive Phosphorus
$ws.Range("F3").Value = 0.191445517051655
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.357142857142857
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0.0003271673143665
$ws.Range("L3").Value = -0.0002998768472906
$ws.Range("M3").Value = 0.0010870535714285
$ws.Range("N3").Value = 2.97424831242348

# Row 4 - E. coli
$ws.Range("D4").Value = $true
$ws.Range("F4").Value = 0.997047304233703
$ws.Range("H4").Value = 0.821428571428571
$ws.Range("J4").Value = 412
$ws.Range("K4").Value = -77.5617626648161
$ws.Range("L4").Value = -133.656688543142
$ws.Range("M4").Value = -40.0400291104939
$ws.Range("N4").Value = -18.8256705497126
$ws.Range("P4").Value = "Virtually certain improving"

# Row 5 - Ammoniacal Nitrogen (NH4)
$ws.Range("F5").Value = 0.9231705779203671
$ws.Range("G5").Value = 0.814814814814815
$ws.Range("H5").Value = 0.222222222222222
$ws.Range("I5").Value = 2
$ws.Range("P5").Value = "Very likely improving"

# Row 6 - Nitrite Nitrogen (NO2)
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = "WARNING: Sen slope influenced by censored values"
$ws.Range("F6").Value = 0.999887112880613
$ws.Range("G6").Value = 0.107142857142857
$ws.Range("J6").Value = 0.0045
$ws.Range("K6").Value = -0.0009487012987012
$ws.Range("L6").Value = -0.001337912087912
$ws.Range("M6").Value = -0.0004933655151879
$ws.Range("N6").Value = -21.0822510822511
$ws.Range("P6").Value = "Virtually certain improving"

# Row 7 - Nitrate Nitrogen (NO3)
$ws.Range("F7").Value = 0.168199467411071
$ws.Range("K7").Value = 0.0415690917314393
$ws.Range("L7").Value = -0.028175637183492
$ws.Range("M7").Value = 0.128944607384299
$ws.Range("N7").Value = 4.51347358647549
$ws.Range("P7").Value = "Unlikely improving"

# Row 8 - pH
$ws.Range("F8").Value = 0.009635406531566
$ws.Range("H8").Value = 0.785714285714286
$ws.Range("J8").Value = 7.125
$ws.Range("K8").Value = -0.07747839731392379
$ws.Range("L8").Value = -0.140966148248435
$ws.Range("M8").Value = -0.0214461758263376
$ws.Range("N8").Value = -1.08741610265156

# Row 9 - SIN (Soluble Inorganic nitrogen)
$ws.Range("F9").Value = 0.168223607188914
$ws.Range("J9").Value = 0.9425
$ws.Range("K9").Value = 0.0422377712330043
$ws.Range("L9").Value = -0.0293515326752144
$ws.Range("M9").Value = 0.127632595646602
$ws.Range("N9").Value = 4.48146113878029
$ws.Range("P9").Value = "Unlikely improving"

# Row 10 - Total Nitrogen
$ws.Range("E10").Value = "ok"
$ws.Range("F10").Value = 0.104084061753305
$ws.Range("H10").Value = 0.821428571428571
$ws.Range("K10").Value = 0.049045466763276
$ws.Range("L10").Value = -0.0129085868713117
$ws.Range("M10").Value = 0.125948275862069
$ws.Range("N10").Value = 4.22805747959276
$ws.Range("P10").Value = "Unlikely improving"

# Row 11 - Total Phosphorus
$ws.Range("F11").Value = 0.6215405335684639
$ws.Range("H11").Value = 0.571428571428571
$ws.Range("K11").Value = -0.0004969387755102
$ws.Range("L11").Value = -0.0025119301487116
$ws.Range("M11").Value = 0.0013719629802319
$ws.Range("N11").Value = -1.65646258503402
$ws.Range("P11").Value = "As likely as not improving"

# Row 12 - ASPM (Macroinvertebrate Average Score Per Metric)
$ws.Range("F12").Value = 0.59675202974633
$ws.Range("K12").Value = 0.0045895084913255
$ws.Range("L12").Value = -0.0237099729999718
$ws.Range("M12").Value = 0.0146490472286005
$ws.Range("N12").Value = 2.92325381613094
$ws.Range("P12").Value = "As likely as not improving"

# Row 13 - MCI (Macroinvertebrate Community Index)
$ws.Range("F13").Value = 0.5
$ws.Range("K13").Value = 0.453163771712159
$ws.Range("M13").Value = 7.06315297117198
$ws.Range("N13").Value = 0.573625027483745
$ws.Range("P13").Value = "As likely as not improving"

# Row 14 - QMCI (Quantitative Macroinvertebrate Community Index)
$ws.Range("F14").Value = 0.04320536648685
$ws.Range("K14").Value = -0.368051321004291
$ws.Range("L14").Value = -1.12283921541135
$ws.Range("M14").Value = -0.0901080576490564
$ws.Range("N14").Value = -8.843136016441409
$ws.Range("P14").Value = "Extremely unlikely improving"
